# Scheduled market-data refresh: update computed price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve
# sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR) with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 154.76923
$ws.Range("I31").Value = 154.76923
$ws.Range("K31").Value = 464.30769
$ws.Range("M31").Value = -234.30769
$ws.Range("H40").Value = 2475.853
$ws.Range("J40").Value = 2830.4167
$ws.Range("L40").Value = 2830.4167
$ws.Range("N40").Value = -3180.4167
$ws.Range("H68").Value = 57925.363
$ws.Range("J68").Value = 57925.363
$ws.Range("L68").Value = 57925.363
$ws.Range("N68").Value = -59423.363
$ws.Range("H71").Value = 57925.363
$ws.Range("J71").Value = 57925.363
$ws.Range("L71").Value = 173776.089
$ws.Range("N71").Value = -181264.089
$ws.Range("H80").Value = 733.3158
$ws.Range("I80").Value = 545.7143
$ws.Range("J80").Value = 842.75
$ws.Range("K80").Value = 1637.1429
$ws.Range("L80").Value = 2528.25
$ws.Range("M80").Value = -639.1428999999998
$ws.Range("N80").Value = -4524.25
$ws.Range("H83").Value = 733.3158
$ws.Range("I83").Value = 545.7143
$ws.Range("J83").Value = 842.75
$ws.Range("K83").Value = 4911.428699999999
$ws.Range("L83").Value = 7584.75
$ws.Range("M83").Value = 80.57130000000052
$ws.Range("N83").Value = -17568.75
$ws.Range("H88").Value = 1003
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 1003
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H125").Value = 2239.4
$ws.Range("J125").Value = 6645.5
$ws.Range("L125").Value = 59809.5
$ws.Range("N125").Value = -64729.5
$ws.Range("H137").Value = 1572.5714
$ws.Range("I137").Value = 1616.7407
$ws.Range("K137").Value = 4850.2221
$ws.Range("M137").Value = -2300.2221
$ws.Range("H138").Value = 3465.65
$ws.Range("I138").Value = 1833.1333
$ws.Range("K138").Value = 5499.3999
$ws.Range("M138").Value = -359.3999000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1646.5405
$ws.Range("I2").Value = 1455.6
$ws.Range("K2").Value = 1455.6
$ws.Range("M2").Value = -1342.6
$ws.Range("H32").Value = 2126.632
$ws.Range("I32").Value = 2068.0532
$ws.Range("K32").Value = 2068.0532
$ws.Range("M32").Value = -1781.0532
$ws.Range("H102").Value = 1002.2143
$ws.Range("I102").Value = 938.88
$ws.Range("K102").Value = 938.88
$ws.Range("M102").Value = 683.12
$ws.Range("H116").Value = 1646.5405
$ws.Range("I116").Value = 1455.6
$ws.Range("K116").Value = 1455.6
$ws.Range("M116").Value = 838.4000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1646.5405
$ws.Range("I3").Value = 1455.6
$ws.Range("K3").Value = 1455.6
$ws.Range("M3").Value = -1341.6
$ws.Range("H134").Value = 5792
$ws.Range("I134").Value = 5067.381
$ws.Range("J134").Value = 7965.857
$ws.Range("K134").Value = 15202.143
$ws.Range("L134").Value = 23897.571
$ws.Range("M134").Value = -12667.143
$ws.Range("N134").Value = -28967.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2054.6667
$ws.Range("I105").Value = 2025.6
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 2025.6
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -278.5999999999999
$ws.Range("N105").Value = -5694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 54985.715
$ws.Range("I122").Value = 73655.87
$ws.Range("K122").Value = 220967.61
$ws.Range("M122").Value = -218517.61
$ws.Range("H132").Value = 3398.3794
$ws.Range("I132").Value = 3354.6956
$ws.Range("J132").Value = 3565.8333
$ws.Range("K132").Value = 10064.0868
$ws.Range("L132").Value = 10697.4999
$ws.Range("M132").Value = -7534.086800000001
$ws.Range("N132").Value = -15757.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2168.2
$ws.Range("J46").Value = 3270
$ws.Range("L46").Value = 3270
$ws.Range("N46").Value = -3646
$ws.Range("H53").Value = 34499.168
$ws.Range("I53").Value = 33498.5
$ws.Range("J53").Value = 36500.5
$ws.Range("K53").Value = 33498.5
$ws.Range("L53").Value = 36500.5
$ws.Range("M53").Value = -32980.5
$ws.Range("N53").Value = -37536.5
$ws.Range("H55").Value = 1526
$ws.Range("I55").Value = 303.85715
$ws.Range("J55").Value = 3427.111
$ws.Range("K55").Value = 303.85715
$ws.Range("L55").Value = 3427.111
$ws.Range("M55").Value = -130.85715
$ws.Range("N55").Value = -3773.111
$ws.Range("H68").Value = 2463.4666
$ws.Range("I68").Value = 723.5714
$ws.Range("J68").Value = 3985.875
$ws.Range("K68").Value = 723.5714
$ws.Range("L68").Value = 3985.875
$ws.Range("M68").Value = 25.42859999999996
$ws.Range("N68").Value = -5483.875
$ws.Range("H71").Value = 2463.4666
$ws.Range("I71").Value = 723.5714
$ws.Range("J71").Value = 3985.875
$ws.Range("K71").Value = 3617.857
$ws.Range("L71").Value = 19929.375
$ws.Range("M71").Value = 126.143
$ws.Range("N71").Value = -27417.375
$ws.Range("H136").Value = 3957.0557
$ws.Range("I136").Value = 3177.457
$ws.Range("K136").Value = 9532.370999999999
$ws.Range("M136").Value = -6982.370999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4467490
$ws.Range("I81").Value = 7145924
$ws.Range("J81").Value = 3433.1667
$ws.Range("K81").Value = 14291848
$ws.Range("L81").Value = 6866.3334
$ws.Range("M81").Value = -14290787
$ws.Range("N81").Value = -8988.3334
$ws.Range("H84").Value = 4467490
$ws.Range("I84").Value = 7145924
$ws.Range("J84").Value = 3433.1667
$ws.Range("K84").Value = 71459240
$ws.Range("L84").Value = 34331.667
$ws.Range("M84").Value = -71453936
$ws.Range("N84").Value = -44939.667
$ws.Range("H86").Value = 69237.88
$ws.Range("J86").Value = 69237.88
$ws.Range("L86").Value = 69237.88
$ws.Range("N86").Value = -71483.88
$ws.Range("H89").Value = 69237.88
$ws.Range("J89").Value = 69237.88
$ws.Range("L89").Value = 346189.4
$ws.Range("N89").Value = -357421.4
$ws.Range("H132").Value = 3850.6667
$ws.Range("I132").Value = 3841
$ws.Range("K132").Value = 11523
$ws.Range("M132").Value = -8993
$ws.Range("H136").Value = 3074.3872
$ws.Range("I136").Value = 2998
$ws.Range("J136").Value = 3213.2727
$ws.Range("K136").Value = 8994
$ws.Range("L136").Value = 9639.8181
$ws.Range("M136").Value = -6444
$ws.Range("N136").Value = -14739.8181
